# Add "Construct Binary Tree from Preorder and Inorder Traversal" to the
# review list (sheet "code nippet"), row 16, following the existing pattern.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("code nippet")
$ws2 = $wb.Worksheets.Item("tips")

# New row 16 on the review-list sheet.
$ws1.Range("A16").Value2 = 15
$ws1.Range("B16").Value2 = "leetcode"
$ws1.Range("C16").Value2 = "Construct Binary Tree from Preorder and Inorder Traversal"

# Column C widened to fit the longer question name.
$ws1.Columns.Item(3).ColumnWidth = 48.25

# Restore the saved cursor positions recorded in the file.
$ws1.Activate() | Out-Null
$ws1.Range("D16").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("B8").Select() | Out-Null

$ws1.Activate() | Out-Null
